# Atualização de bases das ligas, do dia: 14-05-2024 às 01:09
# Australia ALeague - odds/result refresh.
#
# Rows 104/105, 112/113, 124/125 and 159/160 had their whole record
# (every column except the running index in column A) swapped between
# the two rows. Rows 167/168 (upcoming fixtures) got their kickoff
# date/teams/odds refreshed, and the two placeholder match-id strings
# in the shared-string table were renumbered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row28 {
    param($ws, [int]$row, [hashtable]$vals)
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}

# --- Row 104 (was row 105's data) ---------------------------------------
Set-Row28 $ws 104 @{
    B = 7127374
    E = "Central Coast Mariners"
    F = "Western Sydney Wanderers"
    H = 0
    I = "H"
    J = 1.909
    L = 3.6
    M = 2.15
    N = 3.6
    O = 3.25
    P = -0.25
    Q = 1.86
    R = 2.04
    S = 2.75
    T = 1.975
    U = 1.875
    V = 1.15
    X = -1
    Y = 0.8600000000000001
    Z = -1
    AA = -1
    AB = 0.875
}

# --- Row 105 (was row 104's data) ---------------------------------------
Set-Row28 $ws 105 @{
    B = 7127370
    E = "Macarthur FC"
    F = "Wellington Phoenix"
    H = 2
    I = "A"
    J = 2.4
    L = 2.625
    M = 2.375
    N = 3.8
    O = 2.75
    P = 0
    Q = 1.8
    R = 2.05
    S = 3
    T = 1.9
    U = 1.95
    V = -1
    X = 1.75
    Y = -1
    Z = 1.05
    AA = 0
    AB = 0
}

# --- Row 112 (was row 113's data) ---------------------------------------
Set-Row28 $ws 112 @{
    B = 7127379
    E = "Melbourne Victory"
    F = "Central Coast Mariners"
    G = 0
    H = 1
    I = "D"
    K = 3.6
    L = 3.8
    N = 3.6
    O = 4
    Q = 1.9
    R = 1.95
    S = 2.75
    T = 1.925
    U = 1.925
    W = -1
    X = 3
    Z = 0.95
    AA = -1
    AB = 0.925
}

# --- Row 113 (was row 112's data) ---------------------------------------
Set-Row28 $ws 113 @{
    B = 7127376
    E = "Newcastle Jets"
    F = "Macarthur FC"
    G = 2
    H = 2
    I = "A"
    K = 4
    L = 3.4
    N = 4.2
    O = 3.6
    Q = 1.89
    R = 2.01
    S = 3.5
    T = 1.95
    U = 1.9
    W = 3.2
    X = -1
    Z = 1.01
    AA = 0.95
    AB = -1
}

# --- Row 124 (was row 125's data) ---------------------------------------
Set-Row28 $ws 124 @{
    B = 7128012
    E = "Macarthur FC"
    F = "Central Coast Mariners"
    G = 0
    H = 3
    I = "A"
    J = 2.4
    K = 3.5
    L = 2.75
    M = 3.4
    N = 3.75
    O = 2.05
    P = 0.25
    Q = 2.025
    R = 1.825
    S = 3
    T = 2.05
    U = 1.8
    W = -1
    X = 1.05
    Z = 0.825
    AA = 0
    AB = 0
}

# --- Row 125 (was row 124's data) ---------------------------------------
Set-Row28 $ws 125 @{
    B = 7127388
    E = "Sydney FC"
    F = "Brisbane Roar"
    G = 1
    H = 1
    I = "D"
    J = 1.5
    K = 5
    L = 5
    M = 1.533
    N = 5.25
    O = 5
    P = -1
    Q = 1.8
    R = 2.05
    S = 3.5
    T = 1.925
    U = 1.925
    W = 4.25
    X = -1
    Z = 1.05
    AA = -1
    AB = 0.925
}

# --- Row 159 (was row 160's data) ---------------------------------------
Set-Row28 $ws 159 @{
    B = 7127419
    E = "Wellington Phoenix"
    F = "Macarthur FC"
    G = 3
    H = 0
    I = "H"
    J = 1.85
    K = 3.5
    L = 3.9
    M = 1.55
    N = 4.5
    O = 5.25
    P = -1
    Q = 1.89
    R = 2.01
    S = 3.5
    T = 1.9
    U = 1.95
    V = 0.55
    X = -1
    Y = 0.8899999999999999
    Z = -1
    AA = -1
    AB = 0.95
}

# --- Row 160 (was row 159's data) ---------------------------------------
Set-Row28 $ws 160 @{
    B = 7127418
    E = "Newcastle Jets"
    F = "Central Coast Mariners"
    G = 1
    H = 3
    I = "A"
    J = 3.6
    K = 3.25
    L = 2
    M = 4.2
    N = 4
    O = 1.75
    P = 0.75
    Q = 1.85
    R = 2
    S = 3
    T = 1.975
    U = 1.875
    V = -1
    X = 0.75
    Y = -1
    Z = 1
    AA = 0.9750000000000001
    AB = -1
}

# --- Row 167: fixture rescheduled, odds refreshed -----------------------
Set-Row28 $ws 167 @{
    D = 45430.14583333334
    E = "Wellington Phoenix"
    F = "Melbourne Victory"
    J = 2.875
    K = 3.4
    L = 2.375
    M = 2.8
    N = 3.4
    O = 2.5
    P = 0
    Q = 2.03
    R = 1.87
    S = 2.5
    T = 1.85
    U = 2
}

# --- Row 168: fixture rescheduled, odds refreshed -----------------------
Set-Row28 $ws 168 @{
    D = 45430.28125
    E = "Central Coast Mariners"
    F = "Sydney FC"
    K = 4
    L = 3.5
    M = 2.3
    O = 2.8
    P = -0.25
    Q = 2.02
    R = 1.88
    S = 3
    T = 2.05
    U = 1.8
}

# --- Shared-string placeholder match-ids used by rows 167/168 -----------
# These live in the shared-string table as plain text ("t=s" cells) so a
# normal numeric assignment would silently convert them to numbers; force
# Text format for the write, then drop back to the default style so no
# visible formatting change is introduced.
$ws.Range("B167").NumberFormat = "@"
$ws.Range("B167").Value = "8182994"
$ws.Range("B167").Style = "Normal"

$ws.Range("B168").NumberFormat = "@"
$ws.Range("B168").Value = "8182995"
$ws.Range("B168").Style = "Normal"
